$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: manually re-calculated scores -------------------------------
$ws.Range("O7").Value = 862.70809789999998
$ws.Range("P7").Value = 249.45765120005558
# P7 loses its bold/border formatting (matches plain "Normal" cells like P3/P5/P6)
$ws.Range("P7").Style = "Normal"

# --- Row 8: manually re-calculated scores -------------------------------
$ws.Range("F8").Value = 31
$ws.Range("G8").Value = 1058.2903225806451
$ws.Range("H8").Value = 202.21031848851314
$ws.Range("J8").Value = 1281.0147315942022
$ws.Range("P8").Value = 282.63187952757147
# P8 loses its bold/border formatting (matches plain "Normal" cells like P3/P5/P6)
$ws.Range("P8").Style = "Normal"

# --- Column Q: turn the per-row formula into one shared formula --------
# (O2-J2)/P2 filled down Q2:Q8 as a single shared-formula group
$ws.Range("Q2:Q8").Formula = "=(O2-J2)/P2"

# --- Selection moves from E4 to J7 --------------------------------------
$ws.Range("J7").Select()
